# Edit slide 4 ("What we would do differently...") content placeholder:
# trim the three "Headline: long description" bullets down to just their
# bold headlines, and append two more bold headline bullets
# ("More Knowledge" gets its own blank line, then "More Communication").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- Trim "Time management: ..." down to "Time management" ---
$para1 = $tr.Paragraphs(1, 1)
$para1.Runs(2, 1).Text = ""

# --- Trim "Meet in person: ..." down to "Meet in person" ---
$para3 = $tr.Paragraphs(3, 1)
$para3.Runs(2, 1).Text = ""

# --- Trim "More Knowledge: ..." down to "More Knowledge" ---
$para5 = $tr.Paragraphs(5, 1)
$para5.Runs(2, 1).Text = ""

# --- Append a blank (bold) line, "More Communication", then a blank line ---
$lastRun = $para5.Runs(1, 1)
$inserted = $lastRun.InsertAfter("`rBLANK1`rMore Communication`rBLANK2")

# New paragraph: blank bold line right after "More Knowledge"
$para6 = $tr.Paragraphs(6, 1)
$para6.Font.Bold = $true
$para6.Font.Size = 24
$para6.Font.Name = "Cambria"
$para6.Runs(1, 1).Text = ""

# New paragraph: "More Communication" headline, bold
$para7 = $tr.Paragraphs(7, 1)
$para7.Font.Bold = $true
$para7.Font.Size = 24
$para7.Font.Name = "Cambria"

# New paragraph: trailing blank line (not bold), buNone like the other spacer lines
$para8 = $tr.Paragraphs(8, 1)
$para8.Font.Bold = $false
$para8.Font.Size = 24
$para8.Font.Name = "Cambria"
$para8.ParagraphFormat.Bullet.Type = 0
$para8.Runs(1, 1).Text = ""
